{"js": "// Apply the four text replacements described by the diff (v1.4 -> v1.5).\nconst replacements = [\n  {\n    find: \"Administrador esta autenticado no sistema; e, tem permissao para alterar Gerente de Desempenho\",\n    replace: \"Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho\"\n  },\n  {\n    find: \"2. System exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas \",\n    replace: \"2. System exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas \"\n  },\n  {\n    find: \"5. Administrador preenche o campo 'Login do Novo Gerente de Desempenho' do novo Gerente de Desempenho para o Perfil de Competencias  \",\n    replace: \"5. Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias \"\n  },\n  {\n    find: \"2. System apresenta o Catalogo (Perfis) de Competencias cadastradas sem nenhuma alteracao \",\n    replace: \"2. System apresenta o Catalogo (Perfis) de Competencias sem nenhuma alteracao \"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const found = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const range of found.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the four text replacements described by the diff (v1.4 -> v1.5).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Find    = \"Administrador esta autenticado no sistema; e, tem permissao para alterar Gerente de Desempenho\"\n        Replace = \"Administrador esta autenticado no sistema e tem permissao para alterar Gerente de Desempenho\"\n    },\n    @{\n        Find    = \"2. System exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas \"\n        Replace = \"2. System exibe a listagem dos Perfis de Competencias cadastrados com a opcao 'Alterar Gerente' dentre as varias exibidas \"\n    },\n    @{\n        Find    = \"5. Administrador preenche o campo 'Login do Novo Gerente de Desempenho' do novo Gerente de Desempenho para o Perfil de Competencias  \"\n        Replace = \"5. Administrador preenche o campo 'Login do Novo Gerente de Desempenho' para o Perfil de Competencias \"\n    },\n    @{\n        Find    = \"2. System apresenta o Catalogo (Perfis) de Competencias cadastradas sem nenhuma alteracao \"\n        Replace = \"2. System apresenta o Catalogo (Perfis) de Competencias sem nenhuma alteracao \"\n    }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
